# Edit: "Create slice syntax and functions slide"
#
# 1. Slide 4 ("Limitations on Arrays"): merge the split "...10,000 " /
#    "elements over." runs back into a single run, and add a new
#    trailing paragraph "This is where slices come in."
# 2. Insert a new slide 5 ("Slice internals") with a hyperlinked URL.
# 3. Insert a new slide 6 ("Slice Syntax & Functions") with an empty
#    content placeholder.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update slide 4's body text.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2).TextFrame.TextRange

$para1Text = "Arrays are the most primitive data structure and only have a set fixed length. "
$para2Text = "So this leads to a problem, what if we do not know how big of an array we want? We could just create an array that is very large, but that is a waste of computer memory if none of those allocated spots in the array are ever used. We could also create a very small array but what if now we need more size? Copy everything in the current array into another new larger array is costly in time. For example we have an array of 10,000 elements. Now we need to add one more element. We would need to create another array of larger size and copy all 10,000 elements over."
$para3Text = "This is where slices come in."

# Force a real text diff (so the run split collapses into one run) by
# routing through a placeholder value before writing the final text.
$body4.Text = "placeholder"
$body4.Text = $para1Text + "`r" + $para2Text + "`r" + $para3Text

# ---------------------------------------------------------------------
# 2. New slide 5: "Slice internals"
# ---------------------------------------------------------------------
$titleLayout = $p.SlideMaster.CustomLayouts.Item(2)

$s5 = $p.Slides.AddSlide(5, $titleLayout)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Slice internals"

$linkText = "https://blog.golang.org/go-slices-usage-and-internals"
$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = $linkText
$body5.ActionSettings.Item(1).Hyperlink.Address = $linkText

# ---------------------------------------------------------------------
# 3. New slide 6: "Slice Syntax & Functions"
# ---------------------------------------------------------------------
$s6 = $p.Slides.AddSlide(6, $titleLayout)
$title6 = $s6.Shapes.Item(1).TextFrame.TextRange
$title6.Text = "Slice Syntax "
$title6.InsertAfter("& Functions") | Out-Null
